$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('BP1', 'average_doctor_old'),
    @('BQ1', 'average_doctor'),
    @('E4', 0.481),
    @('F4', 0.053),
    @('G4', 0.23),
    @('N4', 0.485),
    @('O4', 0.057),
    @('P4', 0.24),
    @('Q4', 0.052),
    @('R4', 0.035),
    @('S4', 0.188),
    @('W4', 0.367),
    @('AI4', 0.403),
    @('AJ4', 0.096),
    @('AK4', 0.309),
    @('AU4', 0.24),
    @('AW4', 0.163),
    @('BA4', 2.041),
    @('BB4', 0.142),
    @('BC4', 0.377),
    @('BG4', 0.722),
    @('BH4', 0.142),
    @('BI4', 0.377),
    @('BM4', 0.75),
    @('BN4', 0.064),
    @('BO4', 0.253),
    @('BP4', 0.68),
    @('BQ4', 0.761),
    @('E5', 0.606),
    @('F5', 0.058),
    @('G5', 0.241),
    @('N5', 0.733),
    @('O5', 0.065),
    @('P5', 0.255),
    @('Q5', 0.035),
    @('R5', 0.014),
    @('S5', 0.12),
    @('W5', 0.335),
    @('X5', 0.1),
    @('Y5', 0.316),
    @('AI5', 0.404),
    @('AJ5', 0.092),
    @('AK5', 0.303),
    @('AU5', 0.448),
    @('AV5', 0.079),
    @('AW5', 0.28),
    @('BA5', 1.306),
    @('BB5', 0.074),
    @('BC5', 0.271),
    @('BG5', 0.383),
    @('BH5', 0.051),
    @('BI5', 0.225),
    @('BM5', 0.525),
    @('BN5', 0.047),
    @('BO5', 0.217),
    @('BP5', 0.435),
    @('BQ5', 0.458),
    @('E6', 0.536),
    @('N6', 0.584),
    @('Q6', 0.042),
    @('W6', 0.35),
    @('AI6', 0.403),
    @('AU6', 0.313),
    @('BA6', 1.586),
    @('BG6', 0.5),
    @('BM6', 0.618),
    @('BP6', 0.529),
    @('BQ6', 0.569),
    @('E7', 0.576),
    @('N7', 0.665),
    @('Q7', 0.037),
    @('W7', 0.341),
    @('AI7', 0.404),
    @('AU7', 0.382),
    @('BA7', 1.405),
    @('BG7', 0.423),
    @('BM7', 0.559),
    @('BP7', 0.468),
    @('BQ7', 0.497),
    @('E8', 0.705),
    @('F8', 0.072),
    @('G8', 0.268),
    @('N8', 0.821),
    @('O8', 0.044),
    @('P8', 0.21),
    @('Q8', 0.038),
    @('W8', 0.407),
    @('X8', 0.122),
    @('Y8', 0.349),
    @('AI8', 0.472),
    @('AJ8', 0.14),
    @('AK8', 0.374),
    @('AU8', 0.392),
    @('AV8', 0.085),
    @('AW8', 0.292),
    @('BA8', 1.77),
    @('BB8', 0.109),
    @('BC8', 0.33),
    @('BG8', 0.57),
    @('BH8', 0.11),
    @('BI8', 0.331),
    @('BM8', 0.673),
    @('BN8', 0.062),
    @('BO8', 0.249),
    @('BP8', 0.59),
    @('BQ8', 0.626),
    @('E9', 0.667),
    @('F9', 0.222),
    @('G9', 0.471),
    @('N9', 0.762),
    @('O9', 0.181),
    @('P9', 0.426),
    @('W9', 0.31),
    @('X9', 0.214),
    @('Y9', 0.462),
    @('AI9', 0.429),
    @('AJ9', 0.245),
    @('AK9', 0.495),
    @('BA9', 1.738),
    @('BB9', 0.249),
    @('BC9', 0.499),
    @('BG9', 0.619),
    @('BH9', 0.236),
    @('BI9', 0.486),
    @('BM9', 0.643),
    @('BN9', 0.23),
    @('BO9', 0.479),
    @('BP9', 0.579),
    @('BQ9', 0.618),
    @('E10', 0.81),
    @('F10', 0.154),
    @('G10', 0.393),
    @('N10', 0.952),
    @('O10', 0.045),
    @('P10', 0.213),
    @('W10', 0.524),
    @('X10', 0.249),
    @('Y10', 0.499),
    @('AI10', 0.5),
    @('AJ10', 0.25),
    @('AK10', 0.5),
    @('AU10', 0.381),
    @('AV10', 0.236),
    @('AW10', 0.486),
    @('BA10', 2.19),
    @('BB10', 0.214),
    @('BC10', 0.462),
    @('BG10', 0.69),
    @('BH10', 0.214),
    @('BI10', 0.462),
    @('BM10', 0.81),
    @('BN10', 0.154),
    @('BO10', 0.393),
    @('BP10', 0.73),
    @('BQ10', 0.764),
    @('E11', 0.857),
    @('F11', 0.122),
    @('G11', 0.35),
    @('N11', 0.952),
    @('O11', 0.045),
    @('P11', 0.213),
    @('W11', 0.524),
    @('X11', 0.249),
    @('Y11', 0.499),
    @('AI11', 0.571),
    @('AJ11', 0.245),
    @('AK11', 0.495),
    @('AU11', 0.548),
    @('AV11', 0.248),
    @('AW11', 0.498),
    @('BA11', 2.19),
    @('BB11', 0.214),
    @('BC11', 0.462),
    @('BG11', 0.69),
    @('BH11', 0.214),
    @('BI11', 0.462),
    @('BM11', 0.81),
    @('BN11', 0.154),
    @('BO11', 0.393),
    @('BP11', 0.73),
    @('BQ11', 0.77),
    @('E12', 1.417),
    @('F12', 0.854),
    @('G12', 0.924),
    @('N12', 1.25),
    @('O12', 0.287),
    @('P12', 0.536),
    @('W12', 1.5),
    @('X12', 0.432),
    @('Y12', 0.657),
    @('AI12', 1.583),
    @('AJ12', 1.493),
    @('AK12', 1.222),
    @('AU12', 2.88),
    @('AV12', 3.466),
    @('AW12', 1.862),
    @('BA12', 3.787),
    @('BB12', 0.449),
    @('BC12', 0.67),
    @('BG12', 1.138),
    @('BH12', 0.188),
    @('BI12', 0.433),
    @('BM12', 1.235),
    @('BN12', 0.239),
    @('BO12', 0.489),
    @('BP12', 1.262),
    @('BQ12', 1.248),
    @('E13', 1.415),
    @('F13', 0.295),
    @('G13', 0.543),
    @('N13', 1.737),
    @('O13', 0.466),
    @('P13', 0.683),
    @('W13', 0.985),
    @('X13', 0.199),
    @('Y13', 0.446),
    @('AI13', 1.159),
    @('AJ13', 0.312),
    @('AK13', 0.558),
    @('AU13', 2.048),
    @('AV13', 0.344),
    @('AW13', 0.587),
    @('BA13', 2.187),
    @('BB13', 0.278),
    @('BC13', 0.527),
    @('BG13', 0.547),
    @('BH13', 0.05),
    @('BI13', 0.224),
    @('BM13', 0.787),
    @('BN13', 0.163),
    @('BO13', 0.403),
    @('BP13', 0.729),
    @('BQ13', 0.667)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
